$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store its value as literal text (matches the
# source workbook, where Price-column values like "1.001" or "328.15"
# are inline strings, not numbers). Briefly flipping NumberFormat to "@"
# (Text) before the assignment stops Excel from auto-coercing a
# numeric-looking string into a real number; resetting the Style back to
# "Normal" afterwards avoids leaving a stray number-format override on
# the cell.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.883.15"
$ws.Range("E2").Value = "  +1.59%  "

Set-TextValue $ws.Range("D3") "1.770.80"
$ws.Range("E3").Value = "  +2.01%  "

Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.34%  "

Set-TextValue $ws.Range("D5") "328.15"
$ws.Range("E5").Value = "  +1.95%  "

$ws.Range("E6").Value = "  -0.30%  "

Set-TextValue $ws.Range("D7") "0.4485"
$ws.Range("E7").Value = "  -2.59%  "

Set-TextValue $ws.Range("D8") "0.3566"
$ws.Range("E8").Value = "  +1.31%  "

Set-TextValue $ws.Range("D9") "0.07453"
$ws.Range("E9").Value = "  +1.39%  "

Set-TextValue $ws.Range("D10") "42.06"
$ws.Range("E10").Value = "  +1.06%  "

Set-TextValue $ws.Range("D11") "1.099"
$ws.Range("E11").Value = "  +2.07%  "

$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("E13").Value = "  +2.64%  "

Set-TextValue $ws.Range("D14") "6.028"

Set-TextValue $ws.Range("D15") "7.251"
$ws.Range("E15").Value = "  +2.91%  "

Set-TextValue $ws.Range("D16") "1.774.48"
$ws.Range("E16").Value = "  +2.19%  "

$ws.Range("E17").Value = "  +2.65%  "

Set-TextValue $ws.Range("D18") "0.00001060"
$ws.Range("E18").Value = "  +0.81%  "

Set-TextValue $ws.Range("D19") "0.06448"
$ws.Range("E19").Value = "  +1.18%  "

Set-TextValue $ws.Range("D20") "1.000"
$ws.Range("E20").Value = "  -0.33%  "

Set-TextValue $ws.Range("D21") "17.12"
$ws.Range("E21").Value = "  +2.99%  "

Set-TextValue $ws.Range("D22") "5.784"
$ws.Range("E22").Value = "  +1.04%  "

Set-TextValue $ws.Range("D23") "27.929.88"
$ws.Range("E23").Value = "  +1.52%  "

Set-TextValue $ws.Range("D24") "11.31"
$ws.Range("E24").Value = "  +2.10%  "

Set-TextValue $ws.Range("D25") "2.111"
$ws.Range("E25").Value = "  +0.84%  "

Set-TextValue $ws.Range("D26") "163.03"
$ws.Range("E26").Value = "  +0.35%  "

Set-TextValue $ws.Range("D27") "20.37"
$ws.Range("E27").Value = "  +2.68%  "

Set-TextValue $ws.Range("D28") "1.974.95"
$ws.Range("E28").Value = "  +2.08%  "

Set-TextValue $ws.Range("D29") "2.161"
$ws.Range("E29").Value = "  +6.19%  "

Set-TextValue $ws.Range("D30") "125.12"
$ws.Range("E30").Value = "  +0.58%  "

Set-TextValue $ws.Range("D31") "1.107"
$ws.Range("E31").Value = "  +5.96%  "

Set-TextValue $ws.Range("D32") "0.09187"
$ws.Range("E32").Value = "  +0.09%  "

Set-TextValue $ws.Range("D33") "5.621"
$ws.Range("E33").Value = "  +4.15%  "

Set-TextValue $ws.Range("D34") "3.654"
$ws.Range("E34").Value = "  -0.27%  "

Set-TextValue $ws.Range("D35") "11.89"
$ws.Range("E35").Value = "  +2.65%  "

Set-TextValue $ws.Range("D36") "0.02294"
$ws.Range("E36").Value = "  +1.22%  "

Set-TextValue $ws.Range("D37") "0.06093"
$ws.Range("E37").Value = "  +1.91%  "

Set-TextValue $ws.Range("D38") "0.2104"
$ws.Range("E38").Value = "  +1.95%  "

Set-TextValue $ws.Range("D39") "0.6337"
$ws.Range("E39").Value = "  +1.72%  "

Set-TextValue $ws.Range("D40") "4.958"
$ws.Range("E40").Value = "  +1.12%  "

$ws.Range("E41").Value = "  +0.64%  "

Set-TextValue $ws.Range("D42") "1.392"
$ws.Range("E42").Value = "  +1.28%  "

Set-TextValue $ws.Range("D43") "7.912"
$ws.Range("E43").Value = "  +2.52%  "

Set-TextValue $ws.Range("D44") "13.28"
$ws.Range("E44").Value = "  +2.06%  "

Set-TextValue $ws.Range("D45") "3.741"
$ws.Range("E45").Value = "  +1.22%  "

Set-TextValue $ws.Range("D46") "0.5912"
$ws.Range("E46").Value = "  +2.06%  "

Set-TextValue $ws.Range("D47") "122.31"
$ws.Range("E47").Value = "  +0.43%  "

Set-TextValue $ws.Range("D48") "1.959"
$ws.Range("E48").Value = "  +2.07%  "

Set-TextValue $ws.Range("D49") "0.06906"
$ws.Range("E49").Value = "  +1.16%  "

Set-TextValue $ws.Range("D50") "1.138"
$ws.Range("E50").Value = "  +1.79%  "

Set-TextValue $ws.Range("D51") "73.03"
$ws.Range("E51").Value = "  +2.65%  "

